$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 13.98690533333333
$ws.Range("H2").Value = 41.960716
$ws.Range("I2").Value = 0.1880830969574949
$ws.Range("J2").Value = 0.188083096957495
$ws.Range("M2").Value = 0.4890553333333333
$ws.Range("N2").Value = 1.467166
$ws.Range("O2").Value = 0.9644476581758422
$ws.Range("P2").Value = 0.9644476581758422
$ws.Range("Q2").Value = 6.840370650095111
$ws.Range("R2").Value = 61.56333585085599
$ws.Range("S2").Value = 0.1813963024031159
$ws.Range("T2").Value = 0.1813963024031159
$ws.Range("G3").Value = 13.98690533333333
$ws.Range("H3").Value = 41.960716
$ws.Range("I3").Value = 0.1880830969574949
$ws.Range("J3").Value = 0.188083096957495
$ws.Range("O3").Value = 0.03555234182415776
$ws.Range("P3").Value = 0.03555234182415776
$ws.Range("Q3").Value = 0.2521559293493333
$ws.Range("R3").Value = 2.269403364144
$ws.Range("S3").Value = 0.006686794554379067
$ws.Range("T3").Value = 0.006686794554379069
$ws.Range("I4").Value = 0.5894886704907067
$ws.Range("J4").Value = 0.5894886704907067
$ws.Range("M4").Value = 0.4890553333333333
$ws.Range("N4").Value = 1.467166
$ws.Range("O4").Value = 0.9644476581758422
$ws.Range("P4").Value = 0.9644476581758422
$ws.Range("Q4").Value = 21.43903979366889
$ws.Range("R4").Value = 192.95135814302
$ws.Range("S4").Value = 0.5685309677759528
$ws.Range("T4").Value = 0.5685309677759528
$ws.Range("I5").Value = 0.5894886704907067
$ws.Range("J5").Value = 0.5894886704907067
$ws.Range("O5").Value = 0.03555234182415776
$ws.Range("P5").Value = 0.03555234182415776
$ws.Range("S5").Value = 0.02095770271475391
$ws.Range("T5").Value = 0.02095770271475391
$ws.Range("H6").Value = 49.623002
$ws.Range("I6").Value = 0.2224282325517983
$ws.Range("J6").Value = 0.2224282325517984
$ws.Range("M6").Value = 0.4890553333333333
$ws.Range("N6").Value = 1.467166
$ws.Range("O6").Value = 0.9644476581758422
$ws.Range("P6").Value = 0.9644476581758422
$ws.Range("Q6").Value = 8.089464594703555
$ws.Range("R6").Value = 72.805181352332
$ws.Range("S6").Value = 0.2145203879967735
$ws.Range("T6").Value = 0.2145203879967736
$ws.Range("H7").Value = 49.623002
$ws.Range("I7").Value = 0.2224282325517983
$ws.Range("J7").Value = 0.2224282325517984
$ws.Range("O7").Value = 0.03555234182415776
$ws.Range("P7").Value = 0.03555234182415776
$ws.Range("S7").Value = 0.007907844555024789
$ws.Range("T7").Value = 0.007907844555024791
